# Actualización desde MV -datos-
# Appends 4 new daily rows (07-10-2021 .. 13-10-2021) to the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Date = "07-10-2021"; B = 1.95; C = 2.88; D = 3.84; E = 2.4;  F = -0.6  },
    @{ Date = "08-10-2021"; B = 1.92; C = 2.88; D = 3.84; E = 2.4;  F = -1.35 },
    @{ Date = "12-10-2021"; B = 2.14; C = 2.88; D = 3.84; E = 2.4;  F = -1.24 },
    @{ Date = "13-10-2021"; B = 2.07; C = 3.24; D = 3.84; E = 2.4;  F = -0.87 }
)

$startRow = 193
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $cellA = $ws.Cells.Item($r, 1)
    # Force the date-like text to be stored as plain text instead of being
    # auto-converted into a date serial number, then restore the default
    # "Normal" style so no extra formatting is left on the cell.
    $cellA.NumberFormat = "@"
    $cellA.Value = $row.Date
    $cellA.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
}
